# ALTEZZA ALBERO MEDIO.xlsx - avg height math function found
# Rebuilds the Foglio1 sample data with a new data set (rows 2-20), adds a
# second "h/LOG" column (D), a new label in O9, restyles a few B cells,
# updates the chart source ranges, and repositions the chart.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Column A (NN) and Column B (h) raw values for rows 2-20
# ---------------------------------------------------------------------
$A = @{2=1;3=2;4=3;5=6;6=7;7=14;8=30;9=31;10=61;11=123;12=155;13=246;14=495;15=775;16=3876;17=19380;18=96901;19=480456}
$B = @{2=0;3=2;4=1.61;5=3.048;6=3.34;7=4.96;8=6.5519999999999996;9=6.6379999999999999;10=7.61;11=8.0739999999999998;12=8.0500000000000007;13=8.0779999999999994;14=8.0299999999999994;15=8.0730000000000004;16=8.0609999999999999;17=8.07;18=8.0670000000000002;19=8.08;20=8.0679999999999996}

foreach ($r in 2..19) {
    $ws.Cells.Item($r, 1).Value = $A[$r]
}
foreach ($r in 2..20) {
    $ws.Cells.Item($r, 2).Value = $B[$r]
}

# A20 is a formula (=A18*10), not a literal
$ws.Range("A20").Formula = "=A18*10"

# ---------------------------------------------------------------------
# 2. Column C ("c") - rewritten formula B/POWER(LOG(A,2),2), one per row
#    (breaks the old C2:C17 shared-formula group so its si index is free)
# ---------------------------------------------------------------------
foreach ($r in 3..20) {
    $ws.Cells.Item($r, 3).Formula = "=B$r/POWER(LOG(A$r,2),2)"
}

# ---------------------------------------------------------------------
# 3. Column D ("h=OGRANDE(log)") - new column, B/LOG(A,2)
#    D2 is a plain literal 0, D3:D5 individual formulas, D6:D20 filled as
#    one relative formula (creates a shared-formula group).
# ---------------------------------------------------------------------
$ws.Range("D2").Value = 0
$ws.Range("D3").Formula = "=B3/LOG(A3,2)"
$ws.Range("D4").Formula = "=B4/LOG(A4,2)"
$ws.Range("D5").Formula = "=B5/LOG(A5,2)"
$ws.Range("D6:D20").Formula = "=B6/LOG(A6,2)"

# ---------------------------------------------------------------------
# 4. New label cell + restyled cells (#,##0 number format, matches B5/B16
#    in the original file)
# ---------------------------------------------------------------------
$ws.Range("O9").Value = "h=OGRANDE(log)"

$ws.Range("B12").NumberFormat = "#,##0"
$ws.Range("B13").NumberFormat = "#,##0"
$ws.Range("B14").NumberFormat = "#,##0"
$ws.Range("B25").NumberFormat = "#,##0"

# ---------------------------------------------------------------------
# 5. Selection moved to column D
# ---------------------------------------------------------------------
$ws.Columns.Item(4).Select()

# ---------------------------------------------------------------------
# 6. Chart: extend both series to the new data range (A2:A19 / B.. / C..)
# ---------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$ch = $co.Chart
$s1 = $ch.SeriesCollection().Item(1)
$s2 = $ch.SeriesCollection().Item(2)
$s1.Formula = "=SERIES(,Foglio1!`$A`$2:`$A`$19,Foglio1!`$B`$2:`$B`$19,1)"
$s2.Formula = "=SERIES(,Foglio1!`$A`$2:`$A`$19,Foglio1!`$C`$2:`$C`$19,2)"

# ---------------------------------------------------------------------
# 7. Chart: move/resize (new twoCellAnchor: from col5/row11 to col12/row29)
# ---------------------------------------------------------------------
$co.Left = 305.31246062992125
$co.Top = 179.25
$co.Width = 433.0625
$co.Height = 261
